$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly record before the existing row 1195, shifting the
# remaining rows (old 1195-1248) down to (1196-1249) and growing the
# used range from A1:T1248 to A1:T1249.
$ws.Rows("1195:1195").Insert()

$ws.Range("A1195").Value = 5
$ws.Range("B1195").Value = "Macroferia Regional de Talca"
$ws.Range("C1195").Value = "Maule"
$ws.Range("D1195").Value = 44753
$ws.Range("E1195").Value = 7
$ws.Range("F1195").Value = "Fruta"
$ws.Range("G1195").Value = 100102
$ws.Range("H1195").Value = "Cítricos"
$ws.Range("I1195").Value = 100102003
$ws.Range("J1195").Value = "Limón"
$ws.Range("K1195").Value = "Sin especificar"
$ws.Range("L1195").Value = "1a amarillo"
$ws.Range("M1195").Value = 500
$ws.Range("N1195").Value = 4000
$ws.Range("O1195").Value = 4000
$ws.Range("P1195").Value = 4000
$ws.Range("Q1195").Value = "`$/malla 16 kilos"
$ws.Range("R1195").Value = "Provincia de Melipilla"
$ws.Range("S1195").Value = 250
$ws.Range("T1195").Value = 16
